$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$r1 = $ws1.Range("A1")
$text = $r1.Value()

$old1 = "✅ 1000 Bs = 4.01 = 15449.44 pesos"
$new1 = "✅ 1000 Bs = 4.08 = 15701.79 pesos"
$old2 = "✅ 15449.44 pesos = 4.0 = 960.49 Bs"
$new2 = "✅ 15701.79 pesos = 4.05 = 943.17 Bs"

$text = $text.Replace($old1, $new1)
$text = $text.Replace($old2, $new2)

$r1.Value = $text

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 245.195
$ws2.Range("N12").Value = 3879
$ws2.Range("O12").Value = 233.002
